$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16: update I16, add J16
$ws.Range("I16").Value = 0.468739117536998
$ws.Range("J16").Value = 0.2150495036779461

# Row 17: update H17, add I17
$ws.Range("H17").Value = 0.5099036351493167
$ws.Range("I17").Value = 0.24

# Row 18: update G18, add H18
$ws.Range("G18").Value = 0.5604363747513331
$ws.Range("H18").Value = 0.3087982760018804

# Row 19: update F19, add G19
$ws.Range("F19").Value = 0.5999036351493168
$ws.Range("G19").Value = 0.32

# Row 20: update E20, add F20
$ws.Range("E20").Value = 0.6299036351493167
$ws.Range("F20").Value = 0.4476495795507702

# Row 21: update D21, add E21
$ws.Range("D21").Value = 0.3603773643037867
$ws.Range("E21").Value = 0.1088966743764388

# Row 22: update C22, add D22
$ws.Range("C22").Value = 0.4107440146302961
$ws.Range("D22").Value = 0.1461563307127136

# Row 23: update B23, add C23
$ws.Range("B23").Value = 0.25708246933236
$ws.Range("C23").Value = 0.09547648014918764

# Row 24: add B24
$ws.Range("B24").Value = 0.0959495356205764
